$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the
#    "if you make these changes..." paragraph. It is being relocated
#    to the end of the brand new "...repo that you expected" paragraph
#    further down, so drop it from its current spot first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "Pulling from a git repo" heading - the new
#    walkthrough paragraphs about cloning/building from git are
#    inserted directly after it.
# ------------------------------------------------------------------
$insertionPoint = $d.Content
$found = $insertionPoint.Find.Execute("Pulling from a git repo", $true, $false,
                                       $false, $false, $false, $true, 1,
                                       $false, "", 0)
$insertionPoint.Collapse(0)

# ------------------------------------------------------------------
# 3. Insert the new paragraphs (keeping their original run
#    boundaries) plus a trailing blank paragraph, right after the
#    heading. Using a WordOpenXML fragment (instead of typing /
#    InsertAfter calls) preserves every run split and bookmark
#    exactly as authored, since plain text insertion would otherwise
#    silently merge adjacent same-formatted runs together.
# ------------------------------------------------------------------
$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>To create a project based on a git repo it is similar to the test job that we set up, first create an appropriately named freestyle job but this time under source code management select git</w:t></w:r><w:r><w:t xml:space="preserve"> (if you cannot see git in the list of providers you will want to check that the git plugin has been successfully installed and resolve any dependency errors that it might have)</w:t></w:r></w:p><w:p><w:r><w:t>Next input the url of the git project you wish to clone. If your project requires you to be logged in to download it then you create a credentials object using the dropdown below the url box. And specifying the details of the user you want Jenkins to use (it is recommended that you create Jenkins its own user for this purpose)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In the next </w:t></w:r><w:r><w:t>box,</w:t></w:r><w:r><w:t xml:space="preserve"> you can specify what branch you want to </w:t></w:r><w:r><w:t>check out</w:t></w:r><w:r><w:t xml:space="preserve"> when </w:t></w:r><w:r><w:t>the git repository is cloned</w:t></w:r></w:p><w:p><w:r><w:t>To check this is able to pull down the repo corre</w:t></w:r><w:r><w:t>ctly just add a shell script to the build commands that’s runs the command ‘ls’ this way when you run the job you can check in the console that you see the files for the repo that you expected</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xmlFragment)
